$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Employee_code" column (old column C). This shifts old
# Policy_Id/Policy_Status (D/E) left into C/D, leaving 4 used columns.
$ws.Columns("C:C").Delete()

# Insert a new row for "Winston Roy" above the current row 2
# (Vinayaka Naik), pushing every other record down by one.
$ws.Rows("2:2").Insert()

# Populate the new row 2 with the new end-user's data.
$ws.Cells.Item(2, 1).Value = "Winston Roy"
$ws.Cells.Item(2, 2).Value = "pashanwinsty1998@gmail.com"
$ws.Cells.Item(2, 3).Value = "5fd320a7d28d01408a4c4e2a"
$ws.Cells.Item(2, 4).Value = $false

# Replace every employee/policy code value in column C (now "Policy_Id")
# with the shared policy id used across all records.
$ws.Cells.Item(3, 3).Value = "5fd320a7d28d01408a4c4e2a"
$ws.Cells.Item(4, 3).Value = "5fd320a7d28d01408a4c4e2a"
$ws.Cells.Item(5, 3).Value = "5fd320a7d28d01408a4c4e2a"
$ws.Cells.Item(6, 3).Value = "5fd320a7d28d01408a4c4e2a"
$ws.Cells.Item(7, 3).Value = "5fd320a7d28d01408a4c4e2a"
